$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.447.29'
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").Value = '1.941.98'
$ws.Range("E3").Value = '  -0.93%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.616'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.46%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.10'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.79%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.364'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '55.78'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0829'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.103'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.77%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.54'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.99%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '2.231.41'
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.818'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.60'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.22'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.00%  '
$ws.Range("D18").Value = '1.934.72'
$ws.Range("E18").Value = '  -1.27%  '
$ws.Range("D19").Value = '36.370.84'
$ws.Range("E19").Value = '  +1.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.54'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.70%  '
$ws.Range("D21").Value = '0.0₃0859'
$ws.Range("E21").Value = '  +0.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.88'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.92%  '
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.42'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.28'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.72%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '161.69'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.130'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.40'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.98%  '
$ws.Range("E31").Value = '  -1.34%  '
$ws.Range("E32").Value = '  +1.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.65'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0626'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.66%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.24'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.99%  '
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("E37").Value = '  +0.16%  '
$ws.Range("E38").Value = '  -2.91%  '
$ws.Range("E39").Value = '  -5.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0976'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.47%  '
$ws.Range("E42").Value = '  +0.77%  '
$ws.Range("E43").Value = '  -3.88%  '
$ws.Range("E44").Value = '  -0.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '15.95'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("D46").Value = '1.346.61'
$ws.Range("E46").Value = '  +1.17%  '
$ws.Range("E47").Value = '  -5.00%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.46'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.89%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.10'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.89%  '
$ws.Range("E50").Value = '  +0.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '45.19'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.53%  '
